$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells hold text (e.g. "15.60", "1.643.90") where plain
# numeric assignment would coerce to a Double and silently drop formatting
# (trailing zeros, thousand-dot groups). Force Text format first so the
# COM Value setter keeps these as literal strings, matching the source data.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.976.05'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.639.22'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.61'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5088'
$ws.Range('E6').Value = '  +0.68%  '
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2560'
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06344'
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.56'
$ws.Range('E10').Value = '  -0.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07777'
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.273'
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.637.83'
$ws.Range('E13').Value = '  -0.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5421'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.03'
$ws.Range('E15').Value = '  -1.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅7671'
$ws.Range('E16').Value = '  -2.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.997.07'
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.002'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '198.60'
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.410'
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.887'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.031'
$ws.Range('E22').Value = '  +0.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.005'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.878'
$ws.Range('E24').Value = '  +0.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.28'
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1191'
$ws.Range('E26').Value = '  +4.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.806'
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.61'
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.235'
$ws.Range('E29').Value = '  -0.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04894'
$ws.Range('E30').Value = '  -0.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.253'
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.165'
$ws.Range('E32').Value = '  -0.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.520'
$ws.Range('E33').Value = '  -0.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.369'
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9040'
$ws.Range('E35').Value = '  +1.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.588'
$ws.Range('E36').Value = '  -0.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.141.12'
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5435'
$ws.Range('E38').Value = '  -2.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01561'
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.002'
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.527'
$ws.Range('E41').Value = '  -1.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0₈126'
$ws.Range('E42').Value = '  +5.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8093'
$ws.Range('E43').Value = '  -1.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.15'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.405'
$ws.Range('E45').Value = '  -4.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.778.18'
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4529'
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('E48').Value = '  -0.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.86'
$ws.Range('E49').Value = '  -0.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05115'
$ws.Range('E50').Value = '  +1.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.004'
$ws.Range('E51').Value = '  -0.20%  '
